$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newA = @(6097,6047,5968,5957,5848,5824,5811,5780,5649,5596,5626,5650,5580,5525,5577,5655,5565,5634,5707,5687,5821,5852,6012,6037,6316,6460,6647,6829,7063,7369,7385,7492,7599,7736,7782,7792,7798,7769,7738,7696,7635,7527,7453,7452,7406)

$newB = @(45994.95833333334,45994.96875,45994.97916666666,45994.98958333334,45995,45995.01041666666,45995.02083333334,45995.03125,45995.04166666666,45995.05208333334,45995.0625,45995.07291666666,45995.08333333334,45995.09375,45995.10416666666,45995.11458333334,45995.125,45995.13541666666,45995.14583333334,45995.15625,45995.16666666666,45995.17708333334,45995.1875,45995.19791666666,45995.20833333334,45995.21875,45995.22916666666,45995.23958333334,45995.25,45995.26041666666,45995.27083333334,45995.28125,45995.29166666666,45995.30208333334,45995.3125,45995.32291666666,45995.33333333334,45995.34375,45995.35416666666,45995.36458333334,45995.375,45995.38541666666,45995.39583333334,45995.40625,45995.41666666666)

for ($i = 0; $i -lt $newA.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value2 = $newA[$i]
    $ws.Cells.Item($row, 2).Value2 = $newB[$i]
}
